# Generate Report for Handoff
#
# The localization-status report is regenerated: the three "rows" of data
# (keyed by source-file GUID/name) are reordered - the ca3748de... row moves
# to the bottom since it now has a fresh handoff pending, while
# ffff01b68f77... and ffffff2fff4132... shift up - and the ca3748de... row
# gets refreshed status/timestamps reflecting that it is ready to hand off
# again.
#
# Because the workbook's hyperlink relationships (r:id -> external URL) stay
# bound to their original worksheet position, each hyperlink's
# TextToDisplay must be updated to stay in sync with the new cell text, in
# addition to updating the cell values themselves. NOTE: the Hyperlinks
# collection must be (re-)fetched *after* all Range.Value writes on that
# worksheet are done, otherwise the TextToDisplay updates do not stick.

$wb = $excel.ActiveWorkbook

function Set-LinkDisplay {
    param($ws, [string]$cellRef, [string]$text)
    $null = $cellRef -match '^([A-Z]+)(\d+)$'
    $target = '$' + $Matches[1] + '$' + $Matches[2]
    $links = @($ws.Hyperlinks)
    foreach ($h in $links) {
        if ($h.Range.Address() -eq $target) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "ffff01b68f77-91a0-458a-b92c-05df0b7578d0.md"

$wsOverview.Range("A3").Value = "ffffff2fff4132-d4ee-47f0-a36f-9793060cefe8.md"

$wsOverview.Range("A4").Value = "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"

Set-LinkDisplay $wsOverview "A2" "ffff01b68f77-91a0-458a-b92c-05df0b7578d0.md"
Set-LinkDisplay $wsOverview "A3" "ffffff2fff4132-d4ee-47f0-a36f-9793060cefe8.md"
Set-LinkDisplay $wsOverview "A4" "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"

# ---------------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "ffff01b68f77-91a0-458a-b92c-05df0b7578d0.md"
$wsZh.Range("C2").Value = "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.zh-cn.xlf"
$wsZh.Range("D2").Value = "2016-01-27 08:22:47"
$wsZh.Range("E2").Value = "366d26ba-c56f-42f6-8320-c2b4558e46c9.md"
$wsZh.Range("F2").Value = "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-01-27 08:23:50"

$wsZh.Range("A3").Value = "ffffff2fff4132-d4ee-47f0-a36f-9793060cefe8.md"

$wsZh.Range("A4").Value = "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"
$wsZh.Range("B4").Value = "Ready for handoff"
$wsZh.Range("C4").Value = "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.98109d33c3267de9f4d7b2d85aad706e6522b49d.zh-cn.xlf"
$wsZh.Range("D4").Value = "2016-01-27 08:29:56"
$wsZh.Range("E4").Value = "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"
$wsZh.Range("F4").Value = "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.98109d33c3267de9f4d7b2d85aad706e6522b49d.zh-cn.xlf"
$wsZh.Range("G4").Value = "2016-01-27 08:28:59"

Set-LinkDisplay $wsZh "A2" "ffff01b68f77-91a0-458a-b92c-05df0b7578d0.md"
Set-LinkDisplay $wsZh "C2" "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.zh-cn.xlf"
Set-LinkDisplay $wsZh "E2" "366d26ba-c56f-42f6-8320-c2b4558e46c9.md"
Set-LinkDisplay $wsZh "F2" "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.zh-cn.xlf"

Set-LinkDisplay $wsZh "A3" "ffffff2fff4132-d4ee-47f0-a36f-9793060cefe8.md"

Set-LinkDisplay $wsZh "A4" "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"
Set-LinkDisplay $wsZh "C4" "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.98109d33c3267de9f4d7b2d85aad706e6522b49d.zh-cn.xlf"
Set-LinkDisplay $wsZh "E4" "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"
Set-LinkDisplay $wsZh "F4" "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.98109d33c3267de9f4d7b2d85aad706e6522b49d.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "ffff01b68f77-91a0-458a-b92c-05df0b7578d0.md"
$wsDe.Range("C2").Value = "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.de-de.xlf"
$wsDe.Range("D2").Value = "2016-01-27 08:23:00"
$wsDe.Range("E2").Value = "366d26ba-c56f-42f6-8320-c2b4558e46c9.md"
$wsDe.Range("F2").Value = "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.de-de.xlf"
$wsDe.Range("G2").Value = "2016-01-27 08:24:14"

$wsDe.Range("A3").Value = "ffffff2fff4132-d4ee-47f0-a36f-9793060cefe8.md"

$wsDe.Range("A4").Value = "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"
$wsDe.Range("B4").Value = "Ready for handoff"
$wsDe.Range("C4").Value = "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.98109d33c3267de9f4d7b2d85aad706e6522b49d.de-de.xlf"
$wsDe.Range("D4").Value = "2016-01-27 08:30:11"
$wsDe.Range("E4").Value = "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"
$wsDe.Range("F4").Value = "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.98109d33c3267de9f4d7b2d85aad706e6522b49d.de-de.xlf"
$wsDe.Range("G4").Value = "2016-01-27 08:29:21"

Set-LinkDisplay $wsDe "A2" "ffff01b68f77-91a0-458a-b92c-05df0b7578d0.md"
Set-LinkDisplay $wsDe "C2" "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.de-de.xlf"
Set-LinkDisplay $wsDe "E2" "366d26ba-c56f-42f6-8320-c2b4558e46c9.md"
Set-LinkDisplay $wsDe "F2" "366d26ba-c56f-42f6-8320-c2b4558e46c9.f71aeb345eab118bfcd9e86a546e64b28d90b546.de-de.xlf"

Set-LinkDisplay $wsDe "A3" "ffffff2fff4132-d4ee-47f0-a36f-9793060cefe8.md"

Set-LinkDisplay $wsDe "A4" "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"
Set-LinkDisplay $wsDe "C4" "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.98109d33c3267de9f4d7b2d85aad706e6522b49d.de-de.xlf"
Set-LinkDisplay $wsDe "E4" "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.md"
Set-LinkDisplay $wsDe "F4" "ca3748de-aa4b-483d-b1fe-6683ffb5bcf7.98109d33c3267de9f4d7b2d85aad706e6522b49d.de-de.xlf"

$wb.Save()
